# Rename worksheet "Test_Toy" to "Test_OtherSheet".
# This automatically updates all defined names and formulas that
# reference the sheet by name (R_Addin, R_AddinAnotherDef, test_in,
# test_out, testdiagram, and the formulas in column C/E of the sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test_Toy")
$ws.Name = "Test_OtherSheet"

# The two label-cell blocks on this sheet (K3:K5 and K10:K12) contain
# plain text like "Test_Toy!test_in" used as documentation labels next
# to the add-in parameter cells; these are literal strings, not
# formulas, so they must be updated explicitly.
$ws.Range("K3").Value = "Test_OtherSheet!test_in"
$ws.Range("K4").Value = "Test_OtherSheet!test_out"
$ws.Range("K5").Value = "Test_OtherSheet!testdiagram"
$ws.Range("K10").Value = "Test_OtherSheet!test_in"
$ws.Range("K11").Value = "Test_OtherSheet!test_out"
$ws.Range("K12").Value = "Test_OtherSheet!testdiagram"

# Update the sheet's remembered selection/active cell, then restore the
# workbook's originally active sheet so we don't change the active tab.
$ws.Range("M11").Select()
$wb.Worksheets.Item("Main Test").Select()
